$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new task row (row 19) - values inherit column styles (B=17, C=13, E=6)
$ws.Range("B19").Value = "Création d'une class model pour l'import"
$ws.Range("C19").Value = "Une class pour mapper les nom de colonnes dans le fichier csv"
$ws.Range("E19").Value = "21min"

# Match the taller row height used by similar wrapped-text rows (e.g. row 13)
$ws.Rows.Item(19).RowHeight = 30

# Update the selected cell to the new row, as in the authored workbook
$ws.Range("E19").Select()
